$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 166: "Alice Springs" base renamed to "Pine Gap (Alice Springs)" ---
# (location / lat / long are unchanged)
$ws.Cells.Item(166, 2).Value() = "Pine Gap (Alice Springs)"

# --- New bases appended after the current last row (168) ---
$newBases = @(
    @{ Row = 169; Index = 168; Name = "Naval Communication Station Harold E. Holt"; Location = "Exmouth, Western Australia"; Lat = -21.485900000000001; Long = 114.0956; WrapLocation = $true },
    @{ Row = 170; Index = 169; Name = "Al-Harir Air Base"; Location = "Harir, Kurdistan Region, Iraq"; Lat = 36.313699999999997; Long = 44.305; WrapLocation = $false },
    @{ Row = 171; Index = 170; Name = "Aytos Logistics Center"; Location = "Aytos, Bulgaria"; Lat = 42.424100000000003; Long = 27.1358; WrapLocation = $true },
    @{ Row = 172; Index = 171; Name = "Caserma Ederle"; Location = "Vicenza, Italy"; Lat = 45.323270000000001; Long = 11.344262000000001; WrapLocation = $false },
    @{ Row = 173; Index = 172; Name = "Misawa Air Base"; Location = "Misawa, Aomori"; Lat = 40.421900000000001; Long = 141.22190000000001; WrapLocation = $true },
    @{ Row = 174; Index = 173; Name = "Camp Bondsteel"; Location = "Ferizaj, Kosovo"; Lat = 42.2194; Long = 21.149000000000001; WrapLocation = $false },
    @{ Row = 175; Index = 174; Name = "Naval Air Station Keflavik"; Location = "Keflavik International Airport, Iceland"; Lat = 63.590600000000002; Long = 22.361999999999998; WrapLocation = $true }
)

# Source cell whose format (wrap text + vertical-center) is reused for the
# "location" cells that need to wrap, matching the existing style used
# throughout column C (e.g. C167).
$wrapFormatSource = $ws.Cells.Item(167, 3)

foreach ($base in $newBases) {
    $r = $base.Row

    $ws.Cells.Item($r, 1).Value() = $base.Index
    $ws.Cells.Item($r, 2).Value() = $base.Name
    $ws.Cells.Item($r, 3).Value() = $base.Location
    $ws.Cells.Item($r, 4).Value() = $base.Lat
    $ws.Cells.Item($r, 5).Value() = $base.Long

    if ($base.WrapLocation) {
        $wrapFormatSource.Copy()
        $ws.Cells.Item($r, 3).PasteSpecial(-4122)
        $excel.CutCopyMode() = $false
    }
}

# Row 175 wraps to two visual lines (long "location" text), matching the
# row height already used elsewhere in the sheet for two-line wrapped rows.
$ws.Rows.Item(175).RowHeight() = 28.8

# --- Selection / scroll state to match the saved workbook view ---
$ws.Range("A168:A175").Select()
